$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column B ("wiki"), shifting username -> C and cnt -> D
$ws.Columns.Item(2).Insert()

# Column D (cnt) holds numeric-looking strings (e.g. "1092") in the source data -
# force Text format so Excel does not auto-coerce them to numbers.
$ws.Columns.Item(4).NumberFormat = "@"

# Header row
$ws.Range("A1").Value = 'bot'
$ws.Range("B1").Value = 'wiki'
$ws.Range("C1").Value = 'username'
$ws.Range("D1").Value = 'cnt'

# Data rows
$ws.Cells.Item(2, 1).Value = $true
$ws.Cells.Item(2, 2).Value = 'arwiki'
$ws.Cells.Item(2, 3).Value = 'JarBot'
$ws.Cells.Item(2, 4).Value = '1092'
$ws.Cells.Item(3, 1).Value = $true
$ws.Cells.Item(3, 2).Value = 'wikidatawiki'
$ws.Cells.Item(3, 3).Value = 'Research Bot'
$ws.Cells.Item(3, 4).Value = '472'
$ws.Cells.Item(4, 1).Value = $true
$ws.Cells.Item(4, 2).Value = 'wikidatawiki'
$ws.Cells.Item(4, 3).Value = 'BotMultichill'
$ws.Cells.Item(4, 4).Value = '411'
$ws.Cells.Item(5, 1).Value = $true
$ws.Cells.Item(5, 2).Value = 'wikidatawiki'
$ws.Cells.Item(5, 3).Value = 'SuccuBot'
$ws.Cells.Item(5, 4).Value = '333'
$ws.Cells.Item(6, 1).Value = $true
$ws.Cells.Item(6, 2).Value = 'wikidatawiki'
$ws.Cells.Item(6, 3).Value = 'Edoderoobot'
$ws.Cells.Item(6, 4).Value = '257'
$ws.Cells.Item(7, 1).Value = $true
$ws.Cells.Item(7, 2).Value = 'wikidatawiki'
$ws.Cells.Item(7, 3).Value = 'Mr.Ibrahembot'
$ws.Cells.Item(7, 4).Value = '255'
$ws.Cells.Item(8, 1).Value = $true
$ws.Cells.Item(8, 2).Value = 'rowiki'
$ws.Cells.Item(8, 3).Value = 'Andrebot'
$ws.Cells.Item(8, 4).Value = '220'
$ws.Cells.Item(9, 1).Value = $true
$ws.Cells.Item(9, 2).Value = 'commonswiki'
$ws.Cells.Item(9, 3).Value = 'SchlurcherBot'
$ws.Cells.Item(9, 4).Value = '216'
$ws.Cells.Item(10, 1).Value = $true
$ws.Cells.Item(10, 2).Value = 'commonswiki'
$ws.Cells.Item(10, 3).Value = 'KolbertBot'
$ws.Cells.Item(10, 4).Value = '207'
$ws.Cells.Item(11, 1).Value = $false
$ws.Cells.Item(11, 2).Value = 'commonswiki'
$ws.Cells.Item(11, 3).Value = 'Thesupermat'
$ws.Cells.Item(11, 4).Value = '193'
$ws.Cells.Item(12, 1).Value = $true
$ws.Cells.Item(12, 2).Value = 'cebwiki'
$ws.Cells.Item(12, 3).Value = 'Lsjbot'
$ws.Cells.Item(12, 4).Value = '192'
$ws.Cells.Item(13, 1).Value = $false
$ws.Cells.Item(13, 2).Value = 'commonswiki'
$ws.Cells.Item(13, 3).Value = 'Sakhalinio'
$ws.Cells.Item(13, 4).Value = '136'
$ws.Cells.Item(14, 1).Value = $false
$ws.Cells.Item(14, 2).Value = 'trwiki'
$ws.Cells.Item(14, 3).Value = 'Teacher0691'
$ws.Cells.Item(14, 4).Value = '73'
$ws.Cells.Item(15, 1).Value = $false
$ws.Cells.Item(15, 2).Value = 'commonswiki'
$ws.Cells.Item(15, 3).Value = 'StellarD'
$ws.Cells.Item(15, 4).Value = '72'
$ws.Cells.Item(16, 1).Value = $false
$ws.Cells.Item(16, 2).Value = 'enwiki'
$ws.Cells.Item(16, 3).Value = 'Molestash'
$ws.Cells.Item(16, 4).Value = '61'
$ws.Cells.Item(17, 1).Value = $false
$ws.Cells.Item(17, 2).Value = 'wikidatawiki'
$ws.Cells.Item(17, 3).Value = 'Drupol'
$ws.Cells.Item(17, 4).Value = '57'
$ws.Cells.Item(18, 1).Value = $false
$ws.Cells.Item(18, 2).Value = 'frwiktionary'
$ws.Cells.Item(18, 3).Value = 'Arpyia'
$ws.Cells.Item(18, 4).Value = '49'
$ws.Cells.Item(19, 1).Value = $false
$ws.Cells.Item(19, 2).Value = 'elwiki'
$ws.Cells.Item(19, 3).Value = 'ΖῷονΠολιτικόν'
$ws.Cells.Item(19, 4).Value = '48'
$ws.Cells.Item(20, 1).Value = $false
$ws.Cells.Item(20, 2).Value = 'commonswiki'
$ws.Cells.Item(20, 3).Value = 'Zinneke'
$ws.Cells.Item(20, 4).Value = '46'
$ws.Cells.Item(21, 1).Value = $false
$ws.Cells.Item(21, 2).Value = 'hrwiki'
$ws.Cells.Item(21, 3).Value = 'Kubura'
$ws.Cells.Item(21, 4).Value = '41'
